$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "to ask (person に)"
$ws.Range("A4").Value = "to ride; to board (～に)"
$ws.Range("A5").Value = "to do; to perform (～を)"
$ws.Range("A19").Value = "good-looking (conjugates like いい)"
$ws.Range("A39").Value = "fond of; to like (～が)"
$ws.Range("A40").Value = "disgusted with; to dislike (～が)"
$ws.Range("A41").Value = "very fond of; to love (～が)"
$ws.Range("A42").Value = "to hate (～が)"
